# Fixing network data cleaning scripts
# - rename header columns to short machine-friendly names
# - normalize "de"/"el" particles in place names to title case ("De"/"El")
# - normalize "GUANAJUATO" to "Guanajuato"
# - tiny floating point rounding change on two percentage cells
# - remove trailing metadata/footer rows (67-71)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row renames
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Small floating-point adjustments
$ws.Range("D15").Value = 0.09523809523809525
$ws.Range("D19").Value = 0.09523809523809525

# Capitalization fixes for place names
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("B22").Value = "Pánuco De Coronado"
$ws.Range("A25").Value = "Estado De México"
$ws.Range("A27").Value = "Guanajuato"
$ws.Range("B31").Value = "Atoyac De Álvarez"
$ws.Range("B35").Value = "Pachuca De Soto"
$ws.Range("B43").Value = "Huejuquilla El Alto"
$ws.Range("B51").Value = "Oaxaca De Juárez"
$ws.Range("B52").Value = "Ocotlán De Morelos"

# Remove trailing footer/metadata rows (now sheet ends at row 65)
$ws.Range("A67:D71").EntireRow.Delete()
